# Scene 16A edit:
#  1. Merge the three split runs that together spell out Mara's line
#     ("...isn't so bad. " + "You" + " ready to go?") into a single run.
#  2. Drop the stray empty run (no <w:t>) that trails "Pro: Five more
#     minutes." at the end of the scene.
#
# Both paragraphs are located by their (stable) leading text rather than
# a hard-coded paragraph index, then rewritten in place via
# Range.InsertXML so the exact resulting run/paragraph markup (including
# xml:space="preserve") can be controlled precisely.

$d = $word.ActiveDocument

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$maraPrefix = "Mara (neutral smiling): Well, I guess having them once in a while"
$proPrefix  = "Pro: Five more minutes."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t.StartsWith($maraPrefix)) {
        $xml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000004D">' +
               '<w:pPr><w:pageBreakBefore w:val="0"/><w:spacing w:after="200" w:lineRule="auto"/>' +
               '<w:rPr><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
               '</w:pPr>' +
               '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' +
               '<w:rPr><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr>' +
               '<w:t xml:space="preserve">Mara (neutral smiling): Well, I guess having them once in a while isn' + [char]0x2019 + 't so bad. You ready to go?</w:t>' +
               '</w:r></w:p>'
        $p.Range.InsertXML($xml)
    }
    elseif ($t.StartsWith($proPrefix)) {
        $xml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000050">' +
               '<w:pPr><w:pageBreakBefore w:val="0"/><w:spacing w:after="200" w:lineRule="auto"/>' +
               '<w:rPr><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
               '</w:pPr>' +
               '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' +
               '<w:rPr><w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr>' +
               '<w:t xml:space="preserve">Pro: Five more minutes.</w:t>' +
               '</w:r></w:p>'
        $p.Range.InsertXML($xml)
    }
}
